$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 2
    3 = 2
    4 = 2
    5 = 0
    6 = 1
    7 = 3
    8 = 1
    9 = 4
    10 = 1
    11 = 3
    12 = 0
    13 = 2
    14 = 2
    15 = 2
    16 = 0
    17 = 1
    19 = 3
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 3
    26 = 0
    27 = 0
    28 = 3
    29 = 3
    30 = 3
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 6
    36 = 2
    37 = 3
    38 = 4
    39 = 2
    40 = 2
    41 = 2
    42 = 0
    43 = 2
    44 = 1
    45 = 1
    46 = 3
    47 = 2
    48 = 2
    49 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
